# Apply updates to the flaky-tests summary workbook:
#  - RQ1 (Cause of Flakiness) table: bump Memory count, add Concurrency + Async Wait rows
#  - RQ2 (Fix for Flakiness) table: add rows for (Concurrency) locks, (Async Wait) added waitFor,
#       (Memory) added memory for test
#  - Programming Language table: bump Java count, add C++ row
#  - Move the active selection to reflect where the author was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- RQ1: Cause of Flakiness (B3:C16) ---
# Row 5 (Memory) count 1 -> 2
$ws.Range("C5").Value = 2

# Row 8: new Concurrency cause row (shared string append order matters, so write
# the new text values in the same sequence the original author entered them)
$ws.Range("B8").Value = "Concurrency"

# Row 25: new C++ row (Programming Language table)
$ws.Range("B25").Value = "C++"

# Row 12: (Concurrency) locks (RQ2 fix table)
$ws.Range("E12").Value = "(Concurrency) locks"

# Row 9: new Async Wait cause row
$ws.Range("B9").Value = "Async Wait"

# Row 13: (Async Wait) added waitFor
$ws.Range("E13").Value = "(Async Wait) added waitFor"

# Row 14: (Memory) added memory for test
$ws.Range("E14").Value = "(Memory) added memory for test"

# Now fill in the numeric counts for the rows above
$ws.Range("C8").Value = 3
$ws.Range("C25").Value = 4
$ws.Range("F12").Value = 3
$ws.Range("C9").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1

# --- Programming Language (B18:C29) ---
# Row 22 (Java) count 3 -> 4
$ws.Range("C22").Value = 4

# Update the sheet selection/scroll to match where the author ended up
$ws.Range("E16").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
